# Applies the weekly re-shuffle of Fecha/Volumen/Precio/Unidad/Origen values
# across rows 2-37 (column A/B/C/E/F/G/H/I/Q/R remain unchanged).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44165
$ws.Cells.Item(2, 10).Value = 45
$ws.Cells.Item(2, 11).Value = 22000
$ws.Cells.Item(2, 12).Value = 22000
$ws.Cells.Item(2, 13).Value = 22000
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(2, 16).Value = 880

# Row 3
$ws.Cells.Item(3, 4).Value = 44511
$ws.Cells.Item(3, 10).Value = 73
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 17000
$ws.Cells.Item(3, 13).Value = 16479
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(3, 16).Value = 659

# Row 4
$ws.Cells.Item(4, 4).Value = 44410
$ws.Cells.Item(4, 10).Value = 35
$ws.Cells.Item(4, 11).Value = 34000
$ws.Cells.Item(4, 12).Value = 34000
$ws.Cells.Item(4, 13).Value = 34000
$ws.Cells.Item(4, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(4, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(4, 16).Value = 1360

# Row 5
$ws.Cells.Item(5, 4).Value = 44343
$ws.Cells.Item(5, 10).Value = 40
$ws.Cells.Item(5, 11).Value = 28000
$ws.Cells.Item(5, 12).Value = 28000
$ws.Cells.Item(5, 13).Value = 28000
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(5, 16).Value = 1120

# Row 6
$ws.Cells.Item(6, 4).Value = 44160
$ws.Cells.Item(6, 10).Value = 35
$ws.Cells.Item(6, 11).Value = 21000
$ws.Cells.Item(6, 12).Value = 21000
$ws.Cells.Item(6, 13).Value = 21000
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(6, 16).Value = 840

# Row 7
$ws.Cells.Item(7, 4).Value = 44365
$ws.Cells.Item(7, 10).Value = 70
$ws.Cells.Item(7, 11).Value = 22000
$ws.Cells.Item(7, 12).Value = 23000
$ws.Cells.Item(7, 13).Value = 22500
$ws.Cells.Item(7, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(7, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(7, 16).Value = 900

# Row 8
$ws.Cells.Item(8, 4).Value = 44476
$ws.Cells.Item(8, 10).Value = 73
$ws.Cells.Item(8, 11).Value = 23000
$ws.Cells.Item(8, 12).Value = 24000
$ws.Cells.Item(8, 13).Value = 23521
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(8, 16).Value = 941

# Row 9
$ws.Cells.Item(9, 4).Value = 44411
$ws.Cells.Item(9, 10).Value = 35
$ws.Cells.Item(9, 11).Value = 34000
$ws.Cells.Item(9, 12).Value = 34000
$ws.Cells.Item(9, 13).Value = 34000
$ws.Cells.Item(9, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(9, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(9, 16).Value = 1360

# Row 10
$ws.Cells.Item(10, 4).Value = 44253
$ws.Cells.Item(10, 10).Value = 38
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 18000
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(10, 16).Value = 720

# Row 11
$ws.Cells.Item(11, 4).Value = 44370
$ws.Cells.Item(11, 10).Value = 45
$ws.Cells.Item(11, 11).Value = 32000
$ws.Cells.Item(11, 12).Value = 32000
$ws.Cells.Item(11, 13).Value = 32000
$ws.Cells.Item(11, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(11, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(11, 16).Value = 1280

# Row 12
$ws.Cells.Item(12, 4).Value = 44484
$ws.Cells.Item(12, 10).Value = 71
$ws.Cells.Item(12, 11).Value = 29000
$ws.Cells.Item(12, 12).Value = 30000
$ws.Cells.Item(12, 13).Value = 29507
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(12, 16).Value = 1180

# Row 13
$ws.Cells.Item(13, 4).Value = 44412
$ws.Cells.Item(13, 10).Value = 35
$ws.Cells.Item(13, 11).Value = 24000
$ws.Cells.Item(13, 12).Value = 24000
$ws.Cells.Item(13, 13).Value = 24000
$ws.Cells.Item(13, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(13, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(13, 16).Value = 960

# Row 14
$ws.Cells.Item(14, 4).Value = 44473
$ws.Cells.Item(14, 10).Value = 85
$ws.Cells.Item(14, 11).Value = 35000
$ws.Cells.Item(14, 12).Value = 36000
$ws.Cells.Item(14, 13).Value = 35471
$ws.Cells.Item(14, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(14, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(14, 16).Value = 1419

# Row 15
$ws.Cells.Item(15, 4).Value = 44159
$ws.Cells.Item(15, 10).Value = 35
$ws.Cells.Item(15, 11).Value = 22000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 22000
$ws.Cells.Item(15, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(15, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(15, 16).Value = 880

# Row 16
$ws.Cells.Item(16, 4).Value = 44399
$ws.Cells.Item(16, 10).Value = 38
$ws.Cells.Item(16, 11).Value = 33000
$ws.Cells.Item(16, 12).Value = 33000
$ws.Cells.Item(16, 13).Value = 33000
$ws.Cells.Item(16, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(16, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(16, 16).Value = 1320

# Row 17
$ws.Cells.Item(17, 4).Value = 44250
$ws.Cells.Item(17, 10).Value = 38
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 18000
$ws.Cells.Item(17, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(17, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(17, 16).Value = 720

# Row 18
$ws.Cells.Item(18, 4).Value = 44475
$ws.Cells.Item(18, 10).Value = 73
$ws.Cells.Item(18, 11).Value = 25000
$ws.Cells.Item(18, 12).Value = 26000
$ws.Cells.Item(18, 13).Value = 25479
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(18, 16).Value = 1019

# Row 19
$ws.Cells.Item(19, 4).Value = 44483
$ws.Cells.Item(19, 10).Value = 55
$ws.Cells.Item(19, 11).Value = 29000
$ws.Cells.Item(19, 12).Value = 30000
$ws.Cells.Item(19, 13).Value = 29455
$ws.Cells.Item(19, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(19, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 16).Value = 1178

# Row 20
$ws.Cells.Item(20, 4).Value = 44162
$ws.Cells.Item(20, 10).Value = 35
$ws.Cells.Item(20, 11).Value = 17000
$ws.Cells.Item(20, 12).Value = 17000
$ws.Cells.Item(20, 13).Value = 17000
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(20, 16).Value = 680

# Row 21
$ws.Cells.Item(21, 4).Value = 44252
$ws.Cells.Item(21, 10).Value = 40
$ws.Cells.Item(21, 11).Value = 18000
$ws.Cells.Item(21, 12).Value = 19000
$ws.Cells.Item(21, 13).Value = 18625
$ws.Cells.Item(21, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(21, 15).Value = 'Provincia de Talca'
$ws.Cells.Item(21, 16).Value = 745

# Row 22
$ws.Cells.Item(22, 4).Value = 44487
$ws.Cells.Item(22, 10).Value = 73
$ws.Cells.Item(22, 11).Value = 20000
$ws.Cells.Item(22, 12).Value = 21000
$ws.Cells.Item(22, 13).Value = 20521
$ws.Cells.Item(22, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(22, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(22, 16).Value = 821

# Row 23
$ws.Cells.Item(23, 4).Value = 44452
$ws.Cells.Item(23, 10).Value = 70
$ws.Cells.Item(23, 11).Value = 31000
$ws.Cells.Item(23, 12).Value = 32000
$ws.Cells.Item(23, 13).Value = 31500
$ws.Cells.Item(23, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(23, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(23, 16).Value = 1260

# Row 24
$ws.Cells.Item(24, 4).Value = 44509
$ws.Cells.Item(24, 10).Value = 80
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 16000
$ws.Cells.Item(24, 13).Value = 15500
$ws.Cells.Item(24, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(24, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(24, 16).Value = 620

# Row 25
$ws.Cells.Item(25, 4).Value = 44469
$ws.Cells.Item(25, 10).Value = 73
$ws.Cells.Item(25, 11).Value = 28000
$ws.Cells.Item(25, 12).Value = 29000
$ws.Cells.Item(25, 13).Value = 28521
$ws.Cells.Item(25, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(25, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(25, 16).Value = 1141

# Row 26
$ws.Cells.Item(26, 4).Value = 44453
$ws.Cells.Item(26, 10).Value = 73
$ws.Cells.Item(26, 11).Value = 21000
$ws.Cells.Item(26, 12).Value = 22000
$ws.Cells.Item(26, 13).Value = 21521
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 16).Value = 861

# Row 27
$ws.Cells.Item(27, 4).Value = 44481
$ws.Cells.Item(27, 10).Value = 63
$ws.Cells.Item(27, 11).Value = 22000
$ws.Cells.Item(27, 12).Value = 23000
$ws.Cells.Item(27, 13).Value = 22476
$ws.Cells.Item(27, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(27, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(27, 16).Value = 899

# Row 28
$ws.Cells.Item(28, 4).Value = 44406
$ws.Cells.Item(28, 10).Value = 35
$ws.Cells.Item(28, 11).Value = 32000
$ws.Cells.Item(28, 12).Value = 32000
$ws.Cells.Item(28, 13).Value = 32000
$ws.Cells.Item(28, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(28, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(28, 16).Value = 1280

# Row 29
$ws.Cells.Item(29, 4).Value = 44372
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 33000
$ws.Cells.Item(29, 12).Value = 34000
$ws.Cells.Item(29, 13).Value = 33500
$ws.Cells.Item(29, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(29, 16).Value = 1340

# Row 30
$ws.Cells.Item(30, 4).Value = 44376
$ws.Cells.Item(30, 10).Value = 38
$ws.Cells.Item(30, 11).Value = 27000
$ws.Cells.Item(30, 12).Value = 27000
$ws.Cells.Item(30, 13).Value = 27000
$ws.Cells.Item(30, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(30, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(30, 16).Value = 1080

# Row 31
$ws.Cells.Item(31, 4).Value = 44448
$ws.Cells.Item(31, 10).Value = 45
$ws.Cells.Item(31, 11).Value = 32000
$ws.Cells.Item(31, 12).Value = 32000
$ws.Cells.Item(31, 13).Value = 32000
$ws.Cells.Item(31, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(31, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(31, 16).Value = 1280

# Row 32
$ws.Cells.Item(32, 4).Value = 44510
$ws.Cells.Item(32, 10).Value = 73
$ws.Cells.Item(32, 11).Value = 16500
$ws.Cells.Item(32, 12).Value = 17000
$ws.Cells.Item(32, 13).Value = 16740
$ws.Cells.Item(32, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(32, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(32, 16).Value = 670

# Row 33
$ws.Cells.Item(33, 4).Value = 44161
$ws.Cells.Item(33, 10).Value = 35
$ws.Cells.Item(33, 11).Value = 21000
$ws.Cells.Item(33, 12).Value = 21000
$ws.Cells.Item(33, 13).Value = 21000
$ws.Cells.Item(33, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(33, 15).Value = 'Provincia de Quillota'
$ws.Cells.Item(33, 16).Value = 840

# Row 34
$ws.Cells.Item(34, 4).Value = 44468
$ws.Cells.Item(34, 10).Value = 65
$ws.Cells.Item(34, 11).Value = 24000
$ws.Cells.Item(34, 12).Value = 25000
$ws.Cells.Item(34, 13).Value = 24538
$ws.Cells.Item(34, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(34, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(34, 16).Value = 982

# Row 35
$ws.Cells.Item(35, 4).Value = 44515
$ws.Cells.Item(35, 10).Value = 73
$ws.Cells.Item(35, 11).Value = 16000
$ws.Cells.Item(35, 12).Value = 17000
$ws.Cells.Item(35, 13).Value = 16521
$ws.Cells.Item(35, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(35, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(35, 16).Value = 661

# Row 36
$ws.Cells.Item(36, 4).Value = 44181
$ws.Cells.Item(36, 10).Value = 38
$ws.Cells.Item(36, 11).Value = 26000
$ws.Cells.Item(36, 12).Value = 26000
$ws.Cells.Item(36, 13).Value = 26000
$ws.Cells.Item(36, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(36, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(36, 16).Value = 1040

# Row 37
$ws.Cells.Item(37, 4).Value = 44508
$ws.Cells.Item(37, 10).Value = 68
$ws.Cells.Item(37, 11).Value = 16000
$ws.Cells.Item(37, 12).Value = 17000
$ws.Cells.Item(37, 13).Value = 16515
$ws.Cells.Item(37, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(37, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(37, 16).Value = 661
